# DASHBOARD.xlsx - GRADE Slides dashboard update
# - S12 (row 2): priority bump to P0; change summary + source DOI/license
# - S47 (row 7): status/date/summary/source updates
# - Add new rows for S60 and S61 (opening slides for CAC and bempedoic acid blocks)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slides")

# --- Row 2 (S12): priority, change summary, sources ---
$ws.Cells.Item(2, 6).Value = "P0"
$ws.Cells.Item(2, 8).Value = "Cards inferiores em paralelo + rodapé absoluto; placeholder KM (esquemático) pronto para substituir pela figura do paper."
$ws.Cells.Item(2, 9).Value = "Lancet 2025 (SCOT-HEART 10y; 10.1016/S0140-6736(24)02679-5; CC BY 4.0); Circulation 2020 (LAP)"

# --- Row 7 (S47): status, date, change summary, sources ---
$ws.Cells.Item(7, 5).Value = "Atualizado"

# The "Última atualização" column stores dates as plain text (e.g. "2026-01-25"),
# not real Excel dates. Force text formatting before writing so it is not
# auto-converted to a date serial, then restore the standard body-cell
# formatting (copied from an existing text-date cell) so the style matches
# the rest of the column.
$g7 = $ws.Cells.Item(7, 7)
$g7.NumberFormat = "@"
$g7.Value = "2026-01-25"
$ws.Cells.Item(2, 7).Copy()
$g7.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(7, 8).Value = "Reservado espaço inferior + rodapé absoluto (evita corte em projeção)."
$ws.Cells.Item(7, 9).Value = "NEJM 2023 (CLEAR Outcomes; 10.1056/NEJMoa2215024)"

# --- New rows 30 (S60) and 31 (S61) ---
# Copy formatting (styles) from an existing data row so the new rows match
# the established look (A/B/C/D/G/H/I/J = body style, E/F = centered status/priority style)
$ws.Range("A7:J7").Copy()
$ws.Range("A30:J31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(30, 1).Value = "GRADE"
$ws.Cells.Item(30, 2).Value = "S60"
$ws.Cells.Item(30, 3).Value = "Abertura do bloco CAC"
$ws.Cells.Item(30, 4).Value = "CAC / Evidência"
$ws.Cells.Item(30, 5).Value = "Novo"
$ws.Cells.Item(30, 6).Value = "P0"

$g30 = $ws.Cells.Item(30, 7)
$g30.NumberFormat = "@"
$g30.Value = "2026-01-25"
$ws.Cells.Item(2, 7).Copy()
$g30.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(30, 8).Value = "Slide de abertura no padrão (roteiro + mensagem‑chave)."
$ws.Cells.Item(30, 9).Value = "SBC 2025 (Dislipidemias/Prevenção) + MESA/SCOT-HEART (contexto)"
$ws.Cells.Item(30, 10).Value = "-"

$ws.Cells.Item(31, 1).Value = "GRADE"
$ws.Cells.Item(31, 2).Value = "S61"
$ws.Cells.Item(31, 3).Value = "Abertura do bloco ácido bempedóico"
$ws.Cells.Item(31, 4).Value = "SBC 2025"
$ws.Cells.Item(31, 5).Value = "Novo"
$ws.Cells.Item(31, 6).Value = "P0"

$g31 = $ws.Cells.Item(31, 7)
$g31.NumberFormat = "@"
$g31.Value = "2026-01-25"
$ws.Cells.Item(2, 7).Copy()
$g31.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(31, 8).Value = "Slide note + mensagem‑chave (diretriz → RCT → EtD/BR)."
$ws.Cells.Item(31, 9).Value = "SBC 2025 + NEJM 2023 (CLEAR Outcomes)"
$ws.Cells.Item(31, 10).Value = "-"

$wb.Save()
